$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "Merge para ENTREGA FINAL en MASTER":
# Append a "Contacto:" block (heading line, spacer line, a borderless 3-col
# contact table, and a trailing spacer paragraph) right before the document's
# closing section properties.
# ---------------------------------------------------------------------------

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b w:val="0"/>
      <w:bCs/>
      <w:noProof/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b w:val="0"/>
      <w:bCs/>
      <w:noProof/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">  </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b w:val="0"/>
      <w:bCs/>
      <w:noProof/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>Contacto:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b w:val="0"/>
      <w:bCs/>
      <w:noProof/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="Tablaconcuadrcula"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblBorders>
      <w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/>
      <w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/>
    </w:tblBorders>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="3341"/>
    <w:gridCol w:w="3341"/>
    <w:gridCol w:w="3342"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3341" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
          <w:t>Domingo Salvucci</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3341" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="center"/>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
          <w:t xml:space="preserve">Email: </w:t>
        </w:r>
        <w:hyperlink r:id="rIdContactHL" w:history="1">
          <w:r>
            <w:rPr>
              <w:rStyle w:val="Hipervnculo"/>
              <w:rFonts w:cstheme="minorHAnsi"/>
              <w:b w:val="0"/>
              <w:bCs/>
              <w:noProof/>
              <w:sz w:val="24"/>
              <w:szCs w:val="24"/>
            </w:rPr>
            <w:t>salvuccid@gmail.com</w:t>
          </w:r>
        </w:hyperlink>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="3342" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="right"/>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
          <w:t>Telefono</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:rFonts w:cstheme="minorHAnsi"/>
            <w:b w:val="0"/>
            <w:bCs/>
            <w:noProof/>
            <w:sz w:val="24"/>
            <w:szCs w:val="24"/>
          </w:rPr>
          <w:t>: 11 2222 3334444</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b w:val="0"/>
      <w:bCs/>
      <w:noProof/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdContactHL" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="mailto:salvuccid@gmail.com" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target = $d.Range($d.Content.End, $d.Content.End)
$target.InsertXML($xml)

Write-Host "Contacto block inserted."
